$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- demscalar rows: set columns J (10) through AS (45) to 1 ---
$demscalarRows = @(4,5,6,7,8,9,10,12,13,14,15,21,22)
foreach ($r in $demscalarRows) {
    for ($c = 10; $c -le 45; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}

# --- prodinit rows: explicit numeric values per column ---
# row 96
$ws.Cells.Item(96, 10).Value = 2022614.15389004
$ws.Cells.Item(96, 11).Value = 1961197.29578383
$ws.Cells.Item(96, 12).Value = 2266125.4043327
$ws.Cells.Item(96, 13).Value = 1859874.95963218
$ws.Cells.Item(96, 14).Value = 2140274.80148584
$ws.Cells.Item(96, 15).Value = 1940685.00138495
$ws.Cells.Item(96, 16).Value = 1940685
$ws.Cells.Item(96, 17).Value = 2244135.02199179
$ws.Cells.Item(96, 18).Value = 2595033.19545938
$ws.Cells.Item(96, 19).Value = 3000798.62376515
$ws.Cells.Item(96, 20).Value = 3470010.47853522
$ws.Cells.Item(96, 21).Value = 4012589.39063237
$ws.Cells.Item(96, 22).Value = 4640007.20384338
$ws.Cells.Item(96, 23).Value = 5365529.5261411
$ws.Cells.Item(96, 24).Value = 6204496.20682607
$ws.Cells.Item(96, 25).Value = 7174645.66972672
$ws.Cells.Item(96, 26).Value = 8296489.95989327
$ws.Cells.Item(96, 27).Value = 9593748.43346537
$ws.Cells.Item(96, 28).Value = 11093849.2602965
$ws.Cells.Item(96, 29).Value = 12828509.3426955
$ws.Cells.Item(96, 30).Value = 14834404.9116121
$ws.Cells.Item(96, 31).Value = 17153946.9788017
$ws.Cells.Item(96, 32).Value = 19836178.0404955
$ws.Cells.Item(96, 33).Value = 22937808.9917429
$ws.Cells.Item(96, 34).Value = 26524418.175092
$ws.Cells.Item(96, 35).Value = 30671837.9153131
$ws.Cells.Item(96, 36).Value = 35467757.8559167
$ws.Cells.Item(96, 37).Value = 41013578.9971001
$ws.Cells.Item(96, 38).Value = 47426557.6353809
$ws.Cells.Item(96, 39).Value = 54842284.5346207
$ws.Cells.Item(96, 40).Value = 63417551.7459974
$ws.Cells.Item(96, 41).Value = 73333667.6906192
$ws.Cells.Item(96, 42).Value = 84800290.5961691
$ws.Cells.Item(96, 43).Value = 98059861.3386222
$ws.Cells.Item(96, 44).Value = 113392729.413409
$ws.Cells.Item(96, 45).Value = 131123080.415353

# row 97
$ws.Cells.Item(97, 10).Value = 56417.7103924661
$ws.Cells.Item(97, 11).Value = 59928.5924726295
$ws.Cells.Item(97, 12).Value = 62328.1713375157
$ws.Cells.Item(97, 13).Value = 52673.9304765382
$ws.Cells.Item(97, 14).Value = 59796.3725708556
$ws.Cells.Item(97, 15).Value = 67459.6929665381
$ws.Cells.Item(97, 16).Value = 67459.69
$ws.Cells.Item(97, 17).Value = 71454.0035962515
$ws.Cells.Item(97, 18).Value = 75684.8220016001
$ws.Cells.Item(97, 19).Value = 80166.148754671
$ws.Cells.Item(97, 20).Value = 84912.8165488737
$ws.Cells.Item(97, 21).Value = 89940.536326969
$ws.Cells.Item(97, 22).Value = 95265.9492825424
$ws.Cells.Item(97, 23).Value = 100906.681940505
$ws.Cells.Item(97, 24).Value = 106881.404498934
$ws.Cells.Item(97, 25).Value = 113209.892625348
$ws.Cells.Item(97, 26).Value = 119913.092911972
$ws.Cells.Item(97, 27).Value = 127013.19220663
$ws.Cells.Item(97, 28).Value = 134533.691048742
$ws.Cells.Item(97, 29).Value = 142499.4814535
$ws.Cells.Item(97, 30).Value = 150936.929301666
$ws.Cells.Item(97, 31).Value = 159873.961607716
$ws.Cells.Item(97, 32).Value = 169340.158955144
$ws.Cells.Item(97, 33).Value = 179366.85340491
$ws.Cells.Item(97, 34).Value = 189987.232201078
$ws.Cells.Item(97, 35).Value = 201236.447616905
$ws.Cells.Item(97, 36).Value = 213151.733304959
$ws.Cells.Item(97, 37).Value = 225772.527536367
$ws.Cells.Item(97, 38).Value = 239140.603737112
$ws.Cells.Item(97, 39).Value = 253300.208753425
$ws.Cells.Item(97, 40).Value = 268298.209303933
$ws.Cells.Item(97, 41).Value = 284184.247103286
$ws.Cells.Item(97, 42).Value = 301010.903170712
$ws.Cells.Item(97, 43).Value = 318833.871867347
$ws.Cells.Item(97, 44).Value = 337712.145238382
$ws.Cells.Item(97, 45).Value = 357708.208270171

# row 98
$ws.Cells.Item(98, 10).Value = 101881.150923917
$ws.Cells.Item(98, 11).Value = 92683.8919094084
$ws.Cells.Item(98, 12).Value = 89260.1711158246
$ws.Cells.Item(98, 13).Value = 91125.5823798098
$ws.Cells.Item(98, 14).Value = 90115.9977387007
$ws.Cells.Item(98, 15).Value = 92237.8758416919
$ws.Cells.Item(98, 16).Value = 92237.88
$ws.Cells.Item(98, 17).Value = 97297.0315138912
$ws.Cells.Item(98, 18).Value = 102633.672211624
$ws.Cells.Item(98, 19).Value = 108263.022085512
$ws.Cells.Item(98, 20).Value = 114201.13592857
$ws.Cells.Item(98, 21).Value = 120464.949122467
$ws.Cells.Item(98, 22).Value = 127072.325936893
$ws.Cells.Item(98, 23).Value = 134042.110478095
$ws.Cells.Item(98, 24).Value = 141394.180431897
$ws.Cells.Item(98, 25).Value = 149149.503754457
$ws.Cells.Item(98, 26).Value = 157330.198472456
$ws.Cells.Item(98, 27).Value = 165959.595763274
$ws.Cells.Item(98, 28).Value = 175062.306495031
$ws.Cells.Item(98, 29).Value = 184664.291416298
$ws.Cells.Item(98, 30).Value = 194792.935195625
$ws.Cells.Item(98, 31).Value = 205477.124522074
$ws.Cells.Item(98, 32).Value = 216747.330489469
$ws.Cells.Item(98, 33).Value = 228635.695499351
$ws.Cells.Item(98, 34).Value = 241176.12493046
$ws.Cells.Item(98, 35).Value = 254404.383836197
$ws.Cells.Item(98, 36).Value = 268358.19894584
$ws.Cells.Item(98, 37).Value = 283077.366260418
$ws.Cells.Item(98, 38).Value = 298603.864550109
$ws.Cells.Item(98, 39).Value = 314981.975076851
$ws.Cells.Item(98, 40).Value = 332258.407883617
$ws.Cells.Item(98, 41).Value = 350482.435010514
$ws.Cells.Item(98, 42).Value = 369706.031017661
$ws.Cells.Item(98, 43).Value = 389984.021215589
$ws.Cells.Item(98, 44).Value = 411374.23802593
$ws.Cells.Item(98, 45).Value = 433937.685918322

# row 99
$ws.Cells.Item(99, 10).Value = 4664.11736064897
$ws.Cells.Item(99, 11).Value = 4689.74858639489
$ws.Cells.Item(99, 12).Value = 4863.14365503149
$ws.Cells.Item(99, 13).Value = 4541.08735074789
$ws.Cells.Item(99, 14).Value = 4658.74966350976
$ws.Cells.Item(99, 15).Value = 4371.97778857709
$ws.Cells.Item(99, 16).Value = 4371.978
$ws.Cells.Item(99, 17).Value = 4568.08699485427
$ws.Cells.Item(99, 18).Value = 4772.99263458249
$ws.Cells.Item(99, 19).Value = 4987.08950058108
$ws.Cells.Item(99, 20).Value = 5210.78987354892
$ws.Cells.Item(99, 21).Value = 5444.52452740548
$ws.Cells.Item(99, 22).Value = 5688.74355882075
$ws.Cells.Item(99, 23).Value = 5943.91725395463
$ws.Cells.Item(99, 24).Value = 6210.53699407454
$ws.Cells.Item(99, 25).Value = 6489.11620179543
$ws.Cells.Item(99, 26).Value = 6780.1913297642
$ws.Cells.Item(99, 27).Value = 7084.3228936924
$ws.Cells.Item(99, 28).Value = 7402.0965517266
$ws.Cells.Item(99, 29).Value = 7734.12423223491
$ws.Cells.Item(99, 30).Value = 8081.04531218124
$ws.Cells.Item(99, 31).Value = 8443.52784835677
$ws.Cells.Item(99, 32).Value = 8822.2698638393
$ws.Cells.Item(99, 33).Value = 9218.00069215788
$ws.Cells.Item(99, 34).Value = 9631.48238175124
$ws.Cells.Item(99, 35).Value = 10063.5111634244
$ws.Cells.Item(99, 36).Value = 10514.9189836293
$ws.Cells.Item(99, 37).Value = 10986.5751065223
$ws.Cells.Item(99, 38).Value = 11479.3877878832
$ws.Cells.Item(99, 39).Value = 11994.3060241196
$ws.Cells.Item(99, 40).Value = 12532.321379724
$ws.Cells.Item(99, 41).Value = 13094.4698967038
$ws.Cells.Item(99, 42).Value = 13681.8340896599
$ws.Cells.Item(99, 43).Value = 14295.5450303568
$ws.Cells.Item(99, 44).Value = 14936.7845257973
$ws.Cells.Item(99, 45).Value = 15606.7873939977

# row 101
$ws.Cells.Item(101, 10).Value = 304506.612609075
$ws.Cells.Item(101, 11).Value = 298123.515218177
$ws.Cells.Item(101, 12).Value = 357153.828129358
$ws.Cells.Item(101, 13).Value = 344920.677104532
$ws.Cells.Item(101, 14).Value = 376706.052215824
$ws.Cells.Item(101, 15).Value = 400646.343674031
$ws.Cells.Item(101, 16).Value = 400646.3
$ws.Cells.Item(101, 17).Value = 439064.28452445
$ws.Cells.Item(101, 18).Value = 481166.170622235
$ws.Cells.Item(101, 19).Value = 527305.207714688
$ws.Cells.Item(101, 20).Value = 577868.51831969
$ws.Cells.Item(101, 21).Value = 633280.346143815
$ws.Cells.Item(101, 22).Value = 694005.615634115
$ws.Cells.Item(101, 23).Value = 760553.832855422
$ws.Cells.Item(101, 24).Value = 833483.360422882
$ws.Cells.Item(101, 25).Value = 913406.102357883
$ws.Cells.Item(101, 26).Value = 1000992.63817495
$ws.Cells.Item(101, 27).Value = 1096977.84927634
$ws.Cells.Item(101, 28).Value = 1202167.08486185
$ws.Cells.Item(101, 29).Value = 1317442.91908777
$ws.Cells.Item(101, 30).Value = 1443772.55617005
$ws.Cells.Item(101, 31).Value = 1582215.94556306
$ws.Cells.Item(101, 32).Value = 1733934.67530293
$ws.Cells.Item(101, 33).Value = 1900201.71813396
$ws.Cells.Item(101, 34).Value = 2082412.11219128
$ws.Cells.Item(101, 35).Value = 2282094.66585442
$ws.Cells.Item(101, 36).Value = 2500924.78497976
$ws.Cells.Item(101, 37).Value = 2740738.53013644
$ws.Cells.Item(101, 38).Value = 3003548.02179116
$ws.Cells.Item(101, 39).Value = 3291558.32269651
$ws.Cells.Item(101, 40).Value = 3607185.93913196
$ws.Cells.Item(101, 41).Value = 3953079.09622934
$ws.Cells.Item(101, 42).Value = 4332139.95749991
$ws.Cells.Item(101, 43).Value = 4747548.97499235
$ws.Cells.Item(101, 44).Value = 5202791.57438819
$ws.Cells.Item(101, 45).Value = 5701687.39893165

# row 102
$ws.Cells.Item(102, 10).Value = 1537734.41452054
$ws.Cells.Item(102, 11).Value = 861126.541750376
$ws.Cells.Item(102, 12).Value = 718033.440585372
$ws.Cells.Item(102, 13).Value = 2132502.00478701
$ws.Cells.Item(102, 14).Value = 1486029.32223597
$ws.Cells.Item(102, 15).Value = 885789.958655275
$ws.Cells.Item(102, 16).Value = 885790
$ws.Cells.Item(102, 17).Value = 924648.962001985
$ws.Cells.Item(102, 18).Value = 965212.638358243
$ws.Cells.Item(102, 19).Value = 1007555.81364561
$ws.Cells.Item(102, 20).Value = 1051756.55318584
$ws.Cells.Item(102, 21).Value = 1097896.34696945
$ws.Cells.Item(102, 22).Value = 1146060.25989351
$ws.Cells.Item(102, 23).Value = 1196337.08859014
$ws.Cells.Item(102, 24).Value = 1248819.52513502
$ws.Cells.Item(102, 25).Value = 1303604.32793767
$ws.Cells.Item(102, 26).Value = 1360792.50012854
$ws.Cells.Item(102, 27).Value = 1420489.47577184
$ws.Cells.Item(102, 28).Value = 1482805.31424737
$ws.Cells.Item(102, 29).Value = 1547854.90315973
$ws.Cells.Item(102, 30).Value = 1615758.17014905
$ws.Cells.Item(102, 31).Value = 1686640.30399366
$ws.Cells.Item(102, 32).Value = 1760631.9854124
$ws.Cells.Item(102, 33).Value = 1837869.62799204
$ws.Cells.Item(102, 34).Value = 1918495.62968403
$ws.Cells.Item(102, 35).Value = 2002658.6353342
$ws.Cells.Item(102, 36).Value = 2090513.8107295
$ws.Cells.Item(102, 37).Value = 2182223.12866689
$ws.Cells.Item(102, 38).Value = 2277955.66757196
$ws.Cells.Item(102, 39).Value = 2377887.92321764
$ws.Cells.Item(102, 40).Value = 2482204.13411784
$ws.Cells.Item(102, 41).Value = 2591096.62119588
$ws.Cells.Item(102, 42).Value = 2704766.14235386
$ws.Cells.Item(102, 43).Value = 2823422.26259678
$ws.Cells.Item(102, 44).Value = 2947283.74039379
$ws.Cells.Item(102, 45).Value = 3076578.93098866

# row 103
$ws.Cells.Item(103, 10).Value = 349108.005700727
$ws.Cells.Item(103, 11).Value = 368428.288937285
$ws.Cells.Item(103, 12).Value = 378473.780809437
$ws.Cells.Item(103, 13).Value = 344911.216712456
$ws.Cells.Item(103, 14).Value = 291565.779542445
$ws.Cells.Item(103, 15).Value = 271897.181689769
$ws.Cells.Item(103, 16).Value = 271897.2
$ws.Cells.Item(103, 17).Value = 281429.3040633
$ws.Cells.Item(103, 18).Value = 291295.582247825
$ws.Cells.Item(103, 19).Value = 301507.749946374
$ws.Cells.Item(103, 20).Value = 312077.933267056
$ws.Cells.Item(103, 21).Value = 323018.683432049
$ws.Cells.Item(103, 22).Value = 334342.991681139
$ws.Cells.Item(103, 23).Value = 346064.304697749
$ws.Cells.Item(103, 24).Value = 358196.540575766
$ws.Cells.Item(103, 25).Value = 370754.105346136
$ws.Cells.Item(103, 26).Value = 383751.910082835
$ws.Cells.Item(103, 27).Value = 397205.388608542
$ws.Cells.Item(103, 28).Value = 411130.515821033
$ws.Cells.Item(103, 29).Value = 425543.826662057
$ws.Cells.Item(103, 30).Value = 440462.435751219
$ws.Cells.Item(103, 31).Value = 455904.057708176
$ws.Cells.Item(103, 32).Value = 471887.028187295
$ws.Cells.Item(103, 33).Value = 488430.325649728
$ws.Cells.Item(103, 34).Value = 505553.593898775
$ws.Cells.Item(103, 35).Value = 523277.16540528
$ws.Cells.Item(103, 36).Value = 541622.085450767
$ws.Cells.Item(103, 37).Value = 560610.137116979
$ws.Cells.Item(103, 38).Value = 580263.867151491
$ws.Cells.Item(103, 39).Value = 600606.612740122
$ws.Cells.Item(103, 40).Value = 621662.529217912
$ws.Cells.Item(103, 41).Value = 643456.618751601
$ws.Cells.Item(103, 42).Value = 666014.76002764
$ws.Cells.Item(103, 43).Value = 689363.738980999
$ws.Cells.Item(103, 44).Value = 713531.28060126
$ws.Cells.Item(103, 45).Value = 738546.081853758

# row 104
$ws.Cells.Item(104, 10).Value = 82249.41559415
$ws.Cells.Item(104, 11).Value = 91605.9155649208
$ws.Cells.Item(104, 12).Value = 89241.8062948236
$ws.Cells.Item(104, 13).Value = 95008.734448662
$ws.Cells.Item(104, 14).Value = 97151.3442873255
$ws.Cells.Item(104, 15).Value = 104605.146480966
$ws.Cells.Item(104, 16).Value = 104605.1
$ws.Cells.Item(104, 17).Value = 110163.286695482
$ws.Cells.Item(104, 18).Value = 116016.807359784
$ws.Cells.Item(104, 19).Value = 122181.354548396
$ws.Cells.Item(104, 20).Value = 128673.454639948
$ws.Cells.Item(104, 21).Value = 135510.512141365
$ws.Cells.Item(104, 22).Value = 142710.856347165
$ws.Cells.Item(104, 23).Value = 150293.790478002
$ws.Cells.Item(104, 24).Value = 158279.643430185
$ws.Cells.Item(104, 25).Value = 166689.824274898
$ws.Cells.Item(104, 26).Value = 175546.879653239
$ws.Cells.Item(104, 27).Value = 184874.554220942
$ws.Cells.Item(104, 28).Value = 194697.854304819
$ws.Cells.Item(104, 29).Value = 205043.114941594
$ws.Cells.Item(104, 30).Value = 215938.070478834
$ws.Cells.Item(104, 31).Value = 227411.928927261
$ws.Cells.Item(104, 32).Value = 239495.450263767
$ws.Cells.Item(104, 33).Value = 252221.028895063
$ws.Cells.Item(104, 34).Value = 265622.780503018
$ws.Cells.Item(104, 35).Value = 279736.633504533
$ws.Cells.Item(104, 36).Value = 294600.425371122
$ws.Cells.Item(104, 37).Value = 310254.004066436
$ws.Cells.Item(104, 38).Value = 326739.334873653
$ws.Cells.Item(104, 39).Value = 344100.612899154
$ws.Cells.Item(104, 40).Value = 362384.38155406
$ws.Cells.Item(104, 41).Value = 381639.657331285
$ws.Cells.Item(104, 42).Value = 401918.061212616
$ws.Cells.Item(104, 43).Value = 423273.957058095
$ws.Cells.Item(104, 44).Value = 445764.597348714
$ws.Cells.Item(104, 45).Value = 469450.276673149

# row 111
$ws.Cells.Item(111, 10).Value = 1895.23302975242
$ws.Cells.Item(111, 11).Value = 2021.4667617012
$ws.Cells.Item(111, 12).Value = 1743.86959689415
$ws.Cells.Item(111, 13).Value = 1689.20584283519
$ws.Cells.Item(111, 14).Value = 1716.35186816027
$ws.Cells.Item(111, 15).Value = 2122.62937577108
$ws.Cells.Item(111, 16).Value = 2122.629
$ws.Cells.Item(111, 17).Value = 2213.78247421412
$ws.Cells.Item(111, 18).Value = 2308.85041292548
$ws.Cells.Item(111, 19).Value = 2408.00091759627
$ws.Cells.Item(111, 20).Value = 2511.40930858202
$ws.Cells.Item(111, 21).Value = 2619.25843513731
$ws.Cells.Item(111, 22).Value = 2731.73899873434
$ws.Cells.Item(111, 23).Value = 2849.04989026595
$ws.Cells.Item(111, 24).Value = 2971.39854172935
$ws.Cells.Item(111, 25).Value = 3099.0012930125
$ws.Cells.Item(111, 26).Value = 3232.08377443159
$ws.Cells.Item(111, 27).Value = 3370.88130569614
$ws.Cells.Item(111, 28).Value = 3515.63931200701
$ws.Cells.Item(111, 29).Value = 3666.61375802333
$ws.Cells.Item(111, 30).Value = 3824.07160046547
$ws.Cells.Item(111, 31).Value = 3988.29126015446
$ws.Cells.Item(111, 32).Value = 4159.56311432251
$ws.Cells.Item(111, 33).Value = 4338.1900100652
$ws.Cells.Item(111, 34).Value = 4524.48779984308
$ws.Cells.Item(111, 35).Value = 4718.7858999798
$ws.Cells.Item(111, 36).Value = 4921.42787314409
$ws.Cells.Item(111, 37).Value = 5132.77203584577
$ws.Cells.Item(111, 38).Value = 5353.1920920197
$ws.Cells.Item(111, 39).Value = 5583.07779381835
$ws.Cells.Item(111, 40).Value = 5822.83563078102
$ws.Cells.Item(111, 41).Value = 6072.88954859872
$ws.Cells.Item(111, 42).Value = 6333.68169874526
$ws.Cells.Item(111, 43).Value = 6605.67322030037
$ws.Cells.Item(111, 44).Value = 6889.34505534717
$ws.Cells.Item(111, 45).Value = 7185.19879938571

# row 112
$ws.Cells.Item(112, 10).Value = 118898.135334696
$ws.Cells.Item(112, 11).Value = 123532.867350864
$ws.Cells.Item(112, 12).Value = 123368.292862608
$ws.Cells.Item(112, 13).Value = 131057.970202339
$ws.Cells.Item(112, 14).Value = 150752.125286456
$ws.Cells.Item(112, 15).Value = 137443.929905998
$ws.Cells.Item(112, 16).Value = 137443.9
$ws.Cells.Item(112, 17).Value = 141225.911112141
$ws.Cells.Item(112, 18).Value = 145111.990924692
$ws.Cells.Item(112, 19).Value = 149105.003071335
$ws.Cells.Item(112, 20).Value = 153207.889983679
$ws.Cells.Item(112, 21).Value = 157423.675059524
$ws.Cells.Item(112, 22).Value = 161755.464890787
$ws.Cells.Item(112, 23).Value = 166206.451552737
$ws.Cells.Item(112, 24).Value = 170779.91495622
$ws.Cells.Item(112, 25).Value = 175479.225264606
$ws.Cells.Item(112, 26).Value = 180307.845377253
$ws.Cells.Item(112, 27).Value = 185269.333481293
$ws.Cells.Item(112, 28).Value = 190367.34567365
$ws.Cells.Item(112, 29).Value = 195605.638655197
$ws.Cells.Item(112, 30).Value = 200988.072499051
$ws.Cells.Item(112, 31).Value = 206518.613495044
$ws.Cells.Item(112, 32).Value = 212201.337072461
$ws.Cells.Item(112, 33).Value = 218040.430803206
$ws.Cells.Item(112, 34).Value = 224040.197487603
$ws.Cells.Item(112, 35).Value = 230205.058325111
$ws.Cells.Item(112, 36).Value = 236539.55617228
$ws.Cells.Item(112, 37).Value = 243048.358890366
$ws.Cells.Item(112, 38).Value = 249736.262785053
$ws.Cells.Item(112, 39).Value = 256608.196140826
$ws.Cells.Item(112, 40).Value = 263669.222852604
$ws.Cells.Item(112, 41).Value = 270924.546157298
$ws.Cells.Item(112, 42).Value = 278379.512468051
$ws.Cells.Item(112, 43).Value = 286039.615313986
$ws.Cells.Item(112, 44).Value = 293910.499388361
$ws.Cells.Item(112, 45).Value = 301997.964708114

# row 100: numeric J..P, then shared-string "inf" for Q..AS
$ws.Cells.Item(100, 10).Value = 6864.61838711398
$ws.Cells.Item(100, 11).Value = 3596.74954423003
$ws.Cells.Item(100, 12).Value = 7731.38446472809
$ws.Cells.Item(100, 13).Value = 2059.06907675324
$ws.Cells.Item(100, 14).Value = 0
$ws.Cells.Item(100, 15).Value = 582.554496979977
$ws.Cells.Item(100, 16).Value = 582.5545
$ws.Cells.Item(100, 17).Value = "inf"
$ws.Cells.Item(100, 18).Value = "inf"
$ws.Cells.Item(100, 19).Value = "inf"
$ws.Cells.Item(100, 20).Value = "inf"
$ws.Cells.Item(100, 21).Value = "inf"
$ws.Cells.Item(100, 22).Value = "inf"
$ws.Cells.Item(100, 23).Value = "inf"
$ws.Cells.Item(100, 24).Value = "inf"
$ws.Cells.Item(100, 25).Value = "inf"
$ws.Cells.Item(100, 26).Value = "inf"
$ws.Cells.Item(100, 27).Value = "inf"
$ws.Cells.Item(100, 28).Value = "inf"
$ws.Cells.Item(100, 29).Value = "inf"
$ws.Cells.Item(100, 30).Value = "inf"
$ws.Cells.Item(100, 31).Value = "inf"
$ws.Cells.Item(100, 32).Value = "inf"
$ws.Cells.Item(100, 33).Value = "inf"
$ws.Cells.Item(100, 34).Value = "inf"
$ws.Cells.Item(100, 35).Value = "inf"
$ws.Cells.Item(100, 36).Value = "inf"
$ws.Cells.Item(100, 37).Value = "inf"
$ws.Cells.Item(100, 38).Value = "inf"
$ws.Cells.Item(100, 39).Value = "inf"
$ws.Cells.Item(100, 40).Value = "inf"
$ws.Cells.Item(100, 41).Value = "inf"
$ws.Cells.Item(100, 42).Value = "inf"
$ws.Cells.Item(100, 43).Value = "inf"
$ws.Cells.Item(100, 44).Value = "inf"
$ws.Cells.Item(100, 45).Value = "inf"
